# Segunda parte do Trabalho sobre o Inferno das Dependencias
#
# 1) The "_GoBack" bookmark currently sits (collapsed) at the end of the
#    title paragraph ("O que é inferno das dependências?"). It needs to
#    move to the start of the first of the two trailing " " (space-only)
#    paragraphs.
# 2) A large block of new content (about the "Jarbas" / Serenata de Amor
#    case) is inserted right after the "A implementação de conceito..."
#    paragraph and before the trailing " " paragraphs.

$d = $word.ActiveDocument

# --- 1) Remove the bookmark from its current (title) location ---------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2) Insert the new paragraphs after "A implementação de conceito..." ---
# That paragraph is paragraph #7 (1-based) in the original document.
$anchor = $d.Paragraphs(7).Range
$anchor.InsertParagraphAfter()
$p = $d.Paragraphs(8)
$p.Range.Text = "Inferno de Dependências no Caso Jarbas:"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(9)
$p.Range.Text = "Jarbas é um componente de um projeto, desenvolvido por Eduardo Cuducos, chamado Serenata de Amor. Este projeto consiste basicamente em robôs de buscas que utilizando inteligência artificial para investigar, ou como o próprio idealizador sugere, passar um pente fino nos gastos com a chamada verba indenizatória, na qual todos os nossos deputados têm direitos durante o mandato, podendo-a utilizar sem licitação. "

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(10)
$p.Range.Text = "A ideia é utilizar a data science e machine learning, para analisar de forma ágil, fácil e transparente todos os gastos com esta verba, cruzando-os com informações como CNPJ e conteúdo do gasto. Visto que, muitas vezes esta verba é utilizada de forma inadequada, como por exemplo gasto com bebidas alcoólicas e pagamentos de empresas (fechadas) no nome do próprio deputado ou algum parente."

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(11)
$p.Range.Text = "Para este projeto o inferno das dependências ocorreu durante a atualização automática de um pacote, o reprint, que foi da versão 0.3.0 para a 0.4.0, ocasionando um erro ao acessar o banco de dados, acesso esse que já funcionará anteriormente sem quaisquer erros."

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(12)

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(13)

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(14)

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(15)
$p.Range.Text = "O projeto é totalmente escrito em código aberto e está disponível no GitHub"

# --- 3) Re-add the "_GoBack" bookmark, collapsed at the start of the ---
#        first trailing " " paragraph (now shifted down by 8 paragraphs).
$spacePara = $d.Paragraphs(16)
$bmRange = $spacePara.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
